# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Row -> new value for the "展览" sheet (sheetId 1, 48 data rows, header in row 1)
$exhibitionUpdates = @{
    5  = 9
    6  = 197
    7  = 4562
    8  = 192
    9  = 120
    11 = 90
    13 = 693
    14 = 181
    15 = 976
    18 = 153
    19 = 67
    20 = 113
    22 = 3493
    23 = 5840
    25 = 30
    27 = 521
    28 = 43
    29 = 3351
    30 = 358
    31 = 24
    32 = 2466
    35 = 124
    36 = 212
    37 = 260
    38 = 349
    39 = 124
    40 = 1010
    41 = 905
    42 = 18
    43 = 20
    45 = 47
    46 = 468
    48 = 549
}

# Row -> new value for the "全部类型" sheet (one extra row vs. "展览")
$allTypesUpdates = @{
    5  = 9
    6  = 197
    7  = 4562
    8  = 192
    9  = 120
    12 = 90
    14 = 693
    15 = 181
    16 = 976
    19 = 153
    20 = 67
    21 = 113
    23 = 3493
    24 = 5840
    26 = 30
    28 = 521
    29 = 43
    30 = 3351
    31 = 358
    32 = 24
    33 = 2466
    36 = 124
    37 = 212
    38 = 260
    39 = 349
    40 = 124
    41 = 1010
    42 = 905
    43 = 18
    44 = 20
    46 = 47
    47 = 468
    49 = 549
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
